$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the "Förändrad" (changed) date column C for every existing data
#    row (2..410) from 45205 (2023-10-06) to 45206 (2023-10-07).
for ($r = 2; $r -le 410; $r++) {
    $ws.Cells.Item($r, 3).Value = 45206
}

# 2) Row 410 gains an explicit custom row height (ht="15" customHeight="1").
$ws.Rows.Item(410).RowHeight = 15

# 3) Append a brand-new row 411 with a new cutting-notice record.
$row = 411

$ws.Cells.Item($row, 1).Value = "A 48036-2023"

$ws.Cells.Item($row, 2).Value = 45204
$ws.Cells.Item($row, 2).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 3).Value = 45206
$ws.Cells.Item($row, 3).NumberFormat = "YYYY-MM-DD"

$ws.Cells.Item($row, 4).Value = "SÖDERMANLANDS LÄN"
$ws.Cells.Item($row, 5).Value = "VINGÅKER"

$ws.Cells.Item($row, 7).Value = 1.2
$ws.Cells.Item($row, 8).Value = 0
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
$ws.Cells.Item($row, 11).Value = 0
$ws.Cells.Item($row, 12).Value = 0
$ws.Cells.Item($row, 13).Value = 0
$ws.Cells.Item($row, 14).Value = 0
$ws.Cells.Item($row, 15).Value = 0
$ws.Cells.Item($row, 16).Value = 0
$ws.Cells.Item($row, 17).Value = 0

# Column R keeps the same (empty, word-wrapped) style as the rest of the sheet.
$ws.Cells.Item($row, 18).Value = ""
$ws.Cells.Item($row, 18).WrapText = $true
